$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand new rows of accelerometer data at the top of the data block
$ws.Rows("2:3").Insert()

# Newly inserted rows inherit the header row's formatting; strip it back to plain data formatting
$ws.Range("A2:C3").ClearFormats()

# Remove what are now the last three (old) data rows so the data block nets to one fewer row
$ws.Rows("22:24").Delete()

$ws.Range("A2").Value = -2.092850303649902
$ws.Range("B2").Value = 0.9900987625122062
$ws.Range("C2").Value = 1.828120517730715

$ws.Range("A3").Value = -2.105730056762695
$ws.Range("B3").Value = 0.8735208511352541
$ws.Range("C3").Value = 2.453470587730407
